# Generate Report for Handback
#
# Applies the "handback" update for the 4def335c-...md source file: it has
# now been handed back (in sync with en-US) in both the zh-cn and de-de
# target-language sheets. Mirrors what the report generator does when it
# re-runs after a handback event: flips the Status column, fills in the
# "Latest Target File" (with a hyperlink back to the source doc, same as the
# existing "Source File Name" link), copies the xliff name into "Latest
# Handback File", and stamps "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c7aad2335d42f856fbe7b66881b41130ab58017/e2e/4def335c-7744-4909-aa0e-8ddb1ab5a402.md"
$hyperlinkColor = 15570276   # BGR packing of RGB FF6495ED, matching the "Latest Target File" hyperlink style already used in column A

# ---- zh-cn sheet, row 2 (4def335c source file) ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"

$wsZh.Range("I2").Value = "4def335c-7744-4909-aa0e-8ddb1ab5a402.md"
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = $hyperlinkColor
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $hyperlinkTarget, "", "", "4def335c-7744-4909-aa0e-8ddb1ab5a402.md")

$wsZh.Range("J2").Value = $wsZh.Range("G2").Value2
$wsZh.Range("K2").Value = "2016-08-17 02:43:38"

# ---- de-de sheet, row 2 (4def335c source file) ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = "4def335c-7744-4909-aa0e-8ddb1ab5a402.md"
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = $hyperlinkColor
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $hyperlinkTarget, "", "", "4def335c-7744-4909-aa0e-8ddb1ab5a402.md")

$wsDe.Range("J2").Value = $wsDe.Range("G2").Value2
$wsDe.Range("K2").Value = "2016-08-17 02:43:45"

# ---- Widen columns that now hold the longer "Handed back: in sync with
#      en-US" status text / file-name-with-hyperlink content. ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 39.15   # E: zh-cn status column
$wsOverview.Columns.Item(6).ColumnWidth = 39.15   # F: de-de status column

$wsZh.Columns.Item(3).ColumnWidth = 39.15    # C: Status
$wsZh.Columns.Item(9).ColumnWidth = 39.15    # I: Latest Target File
$wsZh.Columns.Item(10).ColumnWidth = 39.15   # J: Latest Handback File

$wsDe.Columns.Item(3).ColumnWidth = 39.15    # C: Status
$wsDe.Columns.Item(9).ColumnWidth = 39.15    # I: Latest Target File
$wsDe.Columns.Item(10).ColumnWidth = 39.15   # J: Latest Handback File
